$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (bold, border, alignment) from existing header cell H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-23
$data = @{
    2  = @(9, 9)
    3  = @(9, 9)
    4  = @(7, 8)
    5  = @(5, 6)
    6  = @(7, 8)
    7  = @(8, 8)
    8  = @(9, 9)
    9  = @(5, 7)
    10 = @(6, 6)
    11 = @(5, 6)
    12 = @(9, 9)
    13 = @(6, 7)
    14 = @(7, 7)
    15 = @(6, 7)
    16 = @(7, 8)
    17 = @(7, 8)
    18 = @(9, 9)
    19 = @(8, 9)
    20 = @(7, 7)
    21 = @(9, 9)
    22 = @(9, 9)
    23 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
